{"js": "// Load the paragraphs of the document body so we can address the\n// Title paragraph (index 0) and the intro paragraph (index 1).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Extend the title \"Reading a Research Paper\" -> \"Reading a Research Paper Part 1\"\n//    Insert each token as its own call so each lands in its own run,\n//    matching how Word splits runs on discrete edits.\nconst titlePara = paragraphs.items[0];\ntitlePara.getRange(\"End\").insertText(\" \", Word.InsertLocation.end);\ntitlePara.getRange(\"End\").insertText(\"Part\", Word.InsertLocation.end);\ntitlePara.getRange(\"End\").insertText(\" \", Word.InsertLocation.end);\ntitlePara.getRange(\"End\").insertText(\"1\", Word.InsertLocation.end);\n\n// 2) Collapse the intro paragraph's three runs (\"...Tuan et\", \" \",\n//    \"al.\\u00A0paper...questions:\") into a single run with the full text.\n//    NOTE: the original text joins \"al.\" and \"paper\" with a NON-BREAKING\n//    SPACE (U+00A0), not a regular space -- preserve it exactly.\nconst introPara = paragraphs.items[1];\nintroPara.getRange(\"Whole\").insertText(\n  \"Read the Background and Methods in the summary at the beginning of the Tuan et al.\\u00A0paper. Then answer the following questions:\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Extend the title \"Reading a Research Paper\" -> \"Reading a Research Paper Part 1\"\n#    Insert each token with its own InsertAfter/Collapse pair so each token\n#    lands in its own run, matching how Word splits runs on discrete edits.\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\n$titleRange.MoveEnd(1, -1)          # exclude the paragraph mark\n$titleRange.Collapse(0)             # wdCollapseEnd -> caret at end of \"Paper\"\n\n$titleRange.InsertAfter(\" \")\n$titleRange.Collapse(0)\n\n$titleRange.InsertAfter(\"Part\")\n$titleRange.Collapse(0)\n\n$titleRange.InsertAfter(\" \")\n$titleRange.Collapse(0)\n\n$titleRange.InsertAfter(\"1\")\n$titleRange.Collapse(0)\n\n# 2) Collapse the intro paragraph's three runs (\"...Tuan et\", \" \",\n#    \"al.<nbsp>paper...questions:\") into a single run with the full text.\n#    NOTE: the original text joins \"al.\" and \"paper\" with a NON-BREAKING\n#    SPACE (U+00A0), not a regular space -- preserve it exactly.\n$nbsp = [char]0x00A0\n$introText = \"Read the Background and Methods in the summary at the beginning of the Tuan et al.\" + $nbsp + \"paper. Then answer the following questions:\"\n\n$introPara = $d.Paragraphs.Item(2)\n\n# The final text is character-for-character identical to the paragraph's\n# current (3-run) text, so writing it directly is seen as a no-op and the\n# runs would stay split. Stage a one-character placeholder first so the\n# write is a genuine content change, forcing Word to collapse the\n# paragraph down to a single run, then set the real text.\n$stageRange = $introPara.Range\n$stageRange.MoveEnd(1, -1)\n$stageRange.Text = \"X\"\n\n$introRange = $introPara.Range\n$introRange.MoveEnd(1, -1)          # exclude the paragraph mark\n$introRange.Text = $introText\n"}
